$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (rows 12-13).
# This pushes the existing rows 12-42 down to rows 14-44, keeping all
# of their data/formatting intact (new weekly observation inserted).
$ws.Rows("12:13").Insert()

# Row 12: newest weekly price observation
$ws.Range("A12").Value = 8
$ws.Range("B12").Value = "Terminal La Palmera de La Serena"
$ws.Range("C12").Value = "Coquimbo"
$ws.Range("D12").Value = 45002
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = "Fruta"
$ws.Range("G12").Value = 100107
$ws.Range("H12").Value = "Otros"
$ws.Range("I12").Value = 100107011
$ws.Range("J12").Value = "Tuna"
$ws.Range("K12").Value = "Sin especificar"
$ws.Range("L12").Value = "Especial"
$ws.Range("M12").Value = 360
$ws.Range("N12").Value = 13000
$ws.Range("O12").Value = 14000
$ws.Range("P12").Value = 13500
$ws.Range("Q12").Value = "`$/caja 18 kilos"
$ws.Range("R12").Value = "Provincia de Limarí"
$ws.Range("S12").Value = 750
$ws.Range("T12").Value = 18

# Row 13: newest weekly price observation
$ws.Range("A13").Value = 8
$ws.Range("B13").Value = "Terminal La Palmera de La Serena"
$ws.Range("C13").Value = "Coquimbo"
$ws.Range("D13").Value = 45002
$ws.Range("E13").Value = 4
$ws.Range("F13").Value = "Fruta"
$ws.Range("G13").Value = 100107
$ws.Range("H13").Value = "Otros"
$ws.Range("I13").Value = 100107011
$ws.Range("J13").Value = "Tuna"
$ws.Range("K13").Value = "Sin especificar"
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 240
$ws.Range("N13").Value = 11000
$ws.Range("O13").Value = 12000
$ws.Range("P13").Value = 11500
$ws.Range("Q13").Value = "`$/caja 18 kilos"
$ws.Range("R13").Value = "Provincia de Limarí"
$ws.Range("S13").Value = 639
$ws.Range("T13").Value = 18
